$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 142 (44090 SMS Morgan), shifting
# rows 142-145 down to 143-146. New row becomes row 142: Redbridge poll
# dated 15 Sep 2020.
$ws.Rows.Item(142).Insert()

$ws.Cells.Item(142, 1).Value = 44089
$ws.Cells.Item(142, 2).Value = "Redbridge"
$ws.Cells.Item(142, 3).Value = 54
$ws.Cells.Item(142, 4).Value = 35.8
$ws.Cells.Item(142, 5).Value = 39.1
$ws.Cells.Item(142, 6).Value = 8
$ws.Cells.Item(142, 7).Formula = "#N/A"
$ws.Cells.Item(142, 8).Value = 17.1

# Append a new final row (147): Redbridge poll dated 15 Jun 2021.
# Copy the date formatting from the row above (the MidDate column style)
# before writing the new value.
$ws.Cells.Item(146, 1).Copy($ws.Cells.Item(147, 1))
$ws.Cells.Item(147, 1).Value = 44362
$ws.Cells.Item(147, 2).Value = "Redbridge"
$ws.Cells.Item(147, 3).Value = 52.4
$ws.Cells.Item(147, 4).Value = 41
$ws.Cells.Item(147, 5).Value = 37
$ws.Cells.Item(147, 6).Value = 12
$ws.Cells.Item(147, 7).Formula = "#N/A"
$ws.Cells.Item(147, 8).Value = 10

# Update the selection to mirror the post-edit cursor position (A148).
$ws.Cells.Item(148, 1).Select()
